$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel coercing
# numeric-looking strings (e.g. "0.516", "18.69") into real numbers, and
# without leaving a stray quote-prefix / number-format style behind.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "67.739.55"
Set-TextCell $ws.Range("E2") "  +1.26%  "

Set-TextCell $ws.Range("D3") "2.493.00"
Set-TextCell $ws.Range("E3") "  +1.55%  "

Set-TextCell $ws.Range("E4") "  -0.02%  "

Set-TextCell $ws.Range("D5") "586.45"
Set-TextCell $ws.Range("E5") "  +1.05%  "

Set-TextCell $ws.Range("D6") "176.20"
Set-TextCell $ws.Range("E6") "  +4.65%  "

Set-TextCell $ws.Range("E7") "  -0.03%  "

Set-TextCell $ws.Range("D8") "0.516"

Set-TextCell $ws.Range("D9") "0.141"
Set-TextCell $ws.Range("E9") "  +5.52%  "

Set-TextCell $ws.Range("E10") "  +1.09%  "

Set-TextCell $ws.Range("E11") "  +4.24%  "

Set-TextCell $ws.Range("D12") "4.95"
Set-TextCell $ws.Range("E12") "  +1.54%  "

Set-TextCell $ws.Range("D13") "25.71"
Set-TextCell $ws.Range("E13") "  +2.75%  "

Set-TextCell $ws.Range("D14") "2.913.65"
Set-TextCell $ws.Range("E14") "  +0.35%  "

Set-TextCell $ws.Range("D15") "67.594.91"
Set-TextCell $ws.Range("E15") "  +1.28%  "

Set-TextCell $ws.Range("E16") "  +2.80%  "

Set-TextCell $ws.Range("D17") "2.493.07"
Set-TextCell $ws.Range("E17") "  +1.84%  "

Set-TextCell $ws.Range("E18") "  +1.88%  "

Set-TextCell $ws.Range("D19") "7.46"
Set-TextCell $ws.Range("E19") "  +1.85%  "

Set-TextCell $ws.Range("D20") "352.10"
Set-TextCell $ws.Range("E20") "  +0.84%  "

Set-TextCell $ws.Range("E21") "  +2.38%  "

Set-TextCell $ws.Range("E22") "  +0.10%  "

Set-TextCell $ws.Range("D23") "70.66"
Set-TextCell $ws.Range("E23") "  +3.05%  "

Set-TextCell $ws.Range("D24") "4.25"
Set-TextCell $ws.Range("E24") "  +1.98%  "

Set-TextCell $ws.Range("E25") "  -0.17%  "

Set-TextCell $ws.Range("D26") "9.23"
Set-TextCell $ws.Range("E26") "  +2.04%  "

Set-TextCell $ws.Range("D27") "2.621.33"
Set-TextCell $ws.Range("E27") "  +1.61%  "

Set-TextCell $ws.Range("D28") "0.997"
Set-TextCell $ws.Range("E28") "  -0.63%  "

Set-TextCell $ws.Range("D29") "0.0₃0912"
Set-TextCell $ws.Range("E29") "  +3.04%  "

Set-TextCell $ws.Range("D30") "509.24"
Set-TextCell $ws.Range("E30") "  +0.44%  "

Set-TextCell $ws.Range("E31") "  +4.19%  "

Set-TextCell $ws.Range("E32") "  +3.39%  "

Set-TextCell $ws.Range("E33") "  +1.96%  "

Set-TextCell $ws.Range("E34") "  +0.01%  "

Set-TextCell $ws.Range("D35") "0.122"
Set-TextCell $ws.Range("E35") "  +7.48%  "

Set-TextCell $ws.Range("D36") "161.34"
Set-TextCell $ws.Range("E36") "  +2.02%  "

Set-TextCell $ws.Range("B37") "EthereumClassic"
Set-TextCell $ws.Range("C37") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D37") "18.45"
Set-TextCell $ws.Range("E37") "  +1.99%  "

Set-TextCell $ws.Range("B38") "WhiteBITCoin"
Set-TextCell $ws.Range("C38") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell $ws.Range("D38") "18.69"
Set-TextCell $ws.Range("E38") "  +0.34%  "

Set-TextCell $ws.Range("E39") "  +1.86%  "

Set-TextCell $ws.Range("B40") "Stacks"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D40") "1.76"
Set-TextCell $ws.Range("E40") "  +6.28%  "

Set-TextCell $ws.Range("B41") "USDe"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws.Range("D41") "1.00"
Set-TextCell $ws.Range("E41") "  +0.00%  "

Set-TextCell $ws.Range("D42") "0.330"
Set-TextCell $ws.Range("E42") "  +2.74%  "

Set-TextCell $ws.Range("E43") "  +3.28%  "

Set-TextCell $ws.Range("D44") "2.44"
Set-TextCell $ws.Range("E44") "  +4.54%  "

Set-TextCell $ws.Range("D45") "144.42"
Set-TextCell $ws.Range("E45") "  +2.75%  "

Set-TextCell $ws.Range("E46") "  +3.19%  "

Set-TextCell $ws.Range("E47") "  +4.40%  "

Set-TextCell $ws.Range("D48") "0.515"
Set-TextCell $ws.Range("E48") "  +2.14%  "

Set-TextCell $ws.Range("E49") "  +2.74%  "

Set-TextCell $ws.Range("E50") "  +2.49%  "

Set-TextCell $ws.Range("D51") "0.586"
Set-TextCell $ws.Range("E51") "  +1.30%  "
